# Natmi following Dr Hou advice
# Replace the 2-row LR-pair result table with the new 6-row table that
# includes the "ECs" sending/target cluster in addition to the existing
# "FAPs" and "sCs" clusters.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each inner array is one full data row (columns A..T) to be written
# starting at worksheet row 2.
$rows = @(
    @("ECs", "Vtn", "Tnfrsf11b", "ECs",  2, 0.6666666666666666, 3.708951,          11.126853, 0.1962880377928,    0.1962880377928,    1, 0.3333333333333333, 0.08241233333333334, 0.247237,  0.04727005612861496, 0.04727005612861496, 0.3056633061290001, 2.750969755161,    0.009278546563841349, 0.009278546563841351),
    @("ECs", "Vtn", "Tnfrsf11b", "FAPs", 2, 0.6666666666666666, 3.708951,          11.126853, 0.1962880377928,    0.1962880377928,    3, 1,                  1.661024,             4.983072,  0.9527299438713851,  0.952729943871385,   6.160656625824001,  55.445909632416,   0.1870094912289586,   0.1870094912289586),
    @("FAPs","Vtn", "Tnfrsf11b", "ECs",  3, 1,                   10.28369633333333, 30.851089, 0.5442419095121536, 0.5442419095121537, 1, 0.3333333333333333, 0.08241233333333334, 0.247237,  0.04727005612861496, 0.04727005612861496, 0.8475034101214445, 7.627530691093001, 0.02572634561018408,  0.02572634561018409),
    @("FAPs","Vtn", "Tnfrsf11b", "FAPs", 3, 1,                   10.28369633333333, 30.851089, 0.5442419095121536, 0.5442419095121537, 3, 1,                  1.661024,             4.983072,  0.9527299438713851,  0.952729943871385,   17.08146641837867,  153.733197765408,  0.5185155639019695,   0.5185155639019695),
    @("sCs", "Vtn", "Tnfrsf11b", "ECs",  3, 1,                   4.902803666666666, 14.708411, 0.2594700526950464, 0.2594700526950464, 1, 0.3333333333333333, 0.08241233333333334, 0.247237,  0.04727005612861496, 0.04727005612861496, 0.4040514900452222, 3.636463410407,    0.01226516395458952,  0.01226516395458952),
    @("sCs", "Vtn", "Tnfrsf11b", "FAPs", 3, 1,                   4.902803666666666, 14.708411, 0.2594700526950464, 0.2594700526950464, 3, 1,                  1.661024,             4.983072,  0.9527299438713851,  0.952729943871385,   8.143674557621333,  73.293071018592,   0.2472048887404568,   0.2472048887404568)
)

$columns = @("A","B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T")

$rowIndex = 2
foreach ($rowValues in $rows) {
    for ($colIndex = 0; $colIndex -lt $columns.Length; $colIndex++) {
        $cellRef = "$($columns[$colIndex])$rowIndex"
        $ws.Range($cellRef).Value = $rowValues[$colIndex]
    }
    $rowIndex++
}
